{"js": "// The document contains a single, single-column table whose rows each hold\n// one benchmark statistic. This mirrors the commit (\"Fixed README.md stats\n// and docx preparation for all Renaissance - JDK 21 - Z GC tests\"): several\n// single-value cells are corrected, one stray row (an extra \"0.00001\"\n// sample) is removed, a new row (\"0.01086\") is inserted, and three rows that\n// used to hold a full tab-separated data dump are collapsed down to just\n// their final summary figure.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\n// --- Simple value corrections (row/column indices are 0-based; these use\n//     the original row positions, before any rows are deleted/inserted\n//     below) ---\ntable.getCell(0, 0).value = \"0M\";\ntable.getCell(1, 0).value = \"0M\";\ntable.getCell(2, 0).value = \"0M\";\ntable.getCell(3, 0).value = \"185\";\ntable.getCell(5, 0).value = \"0.00026\";\ntable.getCell(6, 0).value = \"0.00007\";\ntable.getCell(9, 0).value = \"0.00006\";\ntable.getCell(10, 0).value = \"0.00008\";\ntable.getCell(11, 0).value = \"0.00026\";\n\n// Rows 44-46 (0-based 43-45) used to contain a tab-separated dump of\n// per-iteration numbers; collapse each down to just the trailing summary\n// value.\ntable.getCell(43, 0).value = \"100\";\ntable.getCell(44, 0).value = \"0.01\";\ntable.getCell(45, 0).value = \"767\";\nawait context.sync();\n\n// --- Remove the stray extra sample row (old row 8, 0-based index 7, value\n//     \"0.00001\") ---\ntable.rows.load(\"items\");\nawait context.sync();\ntable.rows.items[7].delete();\nawait context.sync();\n\n// --- Insert a new row holding \"0.01086\" right after the (now renumbered)\n//     row that holds \"0.00026\" (old row 12, 0-based index 11 -> after the\n//     deletion above it sits at 0-based index 10) ---\ntable.rows.load(\"items\");\nawait context.sync();\ntable.rows.items[10].insertRows(\"After\", 1, [[\"0.01086\"]]);\nawait context.sync();\n", "ps1": "# The document contains a single, single-column table whose rows each hold\n# one benchmark statistic. This script updates the values called out in the\n# commit (\"Fixed README.md stats and docx preparation for all Renaissance -\n# JDK 21 - Z GC tests\"): several single-value cells are corrected, one stray\n# row (an extra \"0.00001\" sample) is removed, a new row (\"0.01086\") is\n# inserted, and three rows that used to hold a full tab-separated data dump\n# are collapsed down to just their final summary figure.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# --- Simple value corrections (row indices are from the original table,\n#     before any rows are deleted/inserted below) ---\n$t.Rows.Item(1).Cells.Item(1).Range.Text = \"0M\"\n$t.Rows.Item(2).Cells.Item(1).Range.Text = \"0M\"\n$t.Rows.Item(3).Cells.Item(1).Range.Text = \"0M\"\n$t.Rows.Item(4).Cells.Item(1).Range.Text = \"185\"\n$t.Rows.Item(6).Cells.Item(1).Range.Text = \"0.00026\"\n$t.Rows.Item(7).Cells.Item(1).Range.Text = \"0.00007\"\n$t.Rows.Item(10).Cells.Item(1).Range.Text = \"0.00006\"\n$t.Rows.Item(11).Cells.Item(1).Range.Text = \"0.00008\"\n$t.Rows.Item(12).Cells.Item(1).Range.Text = \"0.00026\"\n\n# Rows 44-46 used to contain a tab-separated dump of per-iteration numbers;\n# collapse each down to just the trailing summary value.\n$t.Rows.Item(44).Cells.Item(1).Range.Text = \"100\"\n$t.Rows.Item(45).Cells.Item(1).Range.Text = \"0.01\"\n$t.Rows.Item(46).Cells.Item(1).Range.Text = \"767\"\n\n# --- Remove the stray extra sample row (old row 8, value \"0.00001\") ---\n$t.Rows.Item(8).Delete()\n\n# --- Insert a new row holding \"0.01086\" right after the (now renumbered)\n#     row that holds \"0.00026\" (old row 12) ---\n$newRow = $t.Rows.Add($t.Rows.Item(12))\n$newRow.Cells.Item(1).Range.Text = \"0.01086\"\n"}
